$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new product row at row 67 ("ريست سبورت ") which pushes the
# previously existing rows 67-73 down by one (to 68-74).
# ---------------------------------------------------------------------------

# 1) Insert a blank row at 67, shifting rows 67+ down to 68+.
$ws.Rows("67:67").Insert()

# 2) Copy the formatting (styles) from the row right below (the row that used
#    to be row 67, now shifted to row 68) onto the freshly inserted blank
#    row so it matches the rest of the data-row styling exactly.
$ws.Range("A68:N68").Copy()
$ws.Range("A67:N67").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Row height for the new row matches the standard data-row height.
$ws.Rows("67:67").RowHeight = 25.5

# 4) Fill in the new row's content.
$ws.Range("A67").Value = 64
$ws.Range("B67").Value = "ريست سبورت "
$ws.Range("H67").Value = "0:0"
$ws.Range("L67").Value = 25
$ws.Range("N67").Value = "1:0"

# 5) Re-create the merges for the new row (B:G, H:K, L:M) matching the
#    pattern used by every other data row.
$ws.Range("B67:G67").Merge()
$ws.Range("H67:K67").Merge()
$ws.Range("L67:M67").Merge()

# ---------------------------------------------------------------------------
# Renumber the "م" (sequence) column for the rows that shifted down
# (old rows 67-71, now 68-72) so the sequence stays 64..69.
# ---------------------------------------------------------------------------
$ws.Range("A68").Value = 65
$ws.Range("A69").Value = 66
$ws.Range("A70").Value = 67
$ws.Range("A71").Value = 68
$ws.Range("A72").Value = 69

# ---------------------------------------------------------------------------
# The totals row (was row 72, now row 73) needs its cached total bumped by
# the new row's "سعر البيع" value (25).
# ---------------------------------------------------------------------------
$ws.Range("K73").Value = 4222.4899999999998

# ---------------------------------------------------------------------------
# A handful of rows pick up slightly different (re-rendered) heights in the
# final layout; set them explicitly to match.
# ---------------------------------------------------------------------------
$ws.Rows("72:72").RowHeight = 25.5
$ws.Rows("73:73").RowHeight = 25.5
$ws.Rows("74:74").RowHeight = 17.25

Write-Host "Row inserted and content updated"
